$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 26.666666
$ws.Cells.Item(5, 9).Value = 35
$ws.Cells.Item(5, 11).Value = 35
$ws.Cells.Item(5, 13).Value = 80
$ws.Cells.Item(112, 8).Value = 16357.654
$ws.Cells.Item(112, 10).Value = 18053.404
$ws.Cells.Item(112, 12).Value = 54160.212
$ws.Cells.Item(112, 14).Value = -56376.212
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 85
$ws.Cells.Item(4, 9).Value = 85
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 85
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 31
$ws.Cells.Item(4, 14).ClearContents()
$ws.Cells.Item(32, 8).Value = 3266.5
$ws.Cells.Item(32, 9).Value = 2715.914
$ws.Cells.Item(32, 10).Value = 10581.429
$ws.Cells.Item(32, 11).Value = 2715.914
$ws.Cells.Item(32, 12).Value = 10581.429
$ws.Cells.Item(32, 13).Value = -2428.914
$ws.Cells.Item(32, 14).Value = -11155.429
$ws.Cells.Item(45, 8).Value = 14537543
$ws.Cells.Item(45, 9).Value = 19667682
$ws.Cells.Item(45, 10).Value = 2152
$ws.Cells.Item(45, 11).Value = 19667682
$ws.Cells.Item(45, 12).Value = 2152
$ws.Cells.Item(45, 13).Value = -19667305
$ws.Cells.Item(45, 14).Value = -2906
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 200001800
$ws.Cells.Item(105, 9).Value = 2166.6667
$ws.Cells.Item(105, 10).Value = 500001250
$ws.Cells.Item(105, 11).Value = 2166.6667
$ws.Cells.Item(105, 12).Value = 500001250
$ws.Cells.Item(105, 13).Value = -419.6667000000002
$ws.Cells.Item(105, 14).Value = -500004744
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 8626.666999999999
$ws.Cells.Item(7, 9).Value = 70
$ws.Cells.Item(7, 10).Value = 12905
$ws.Cells.Item(7, 11).Value = 70
$ws.Cells.Item(7, 12).Value = 12905
$ws.Cells.Item(7, 13).Value = 43
$ws.Cells.Item(7, 14).Value = -13131
$ws.Cells.Item(31, 8).Value = 2110.7874
$ws.Cells.Item(31, 9).Value = 2182
$ws.Cells.Item(31, 10).Value = 2102.3096
$ws.Cells.Item(31, 11).Value = 2182
$ws.Cells.Item(31, 12).Value = 2102.3096
$ws.Cells.Item(31, 13).Value = -1887
$ws.Cells.Item(31, 14).Value = -2692.3096
$ws.Cells.Item(34, 8).Value = 2110.7874
$ws.Cells.Item(34, 9).Value = 2182
$ws.Cells.Item(34, 10).Value = 2102.3096
$ws.Cells.Item(34, 11).Value = 2182
$ws.Cells.Item(34, 12).Value = 2102.3096
$ws.Cells.Item(34, 13).Value = -1980
$ws.Cells.Item(34, 14).Value = -2506.3096
$ws.Cells.Item(141, 8).Value = 40061.25
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 40061.25
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 40061.25
$ws.Cells.Item(141, 13).ClearContents()
$ws.Cells.Item(141, 14).Value = -50421.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 106.947365
$ws.Cells.Item(2, 10).Value = 184.2
$ws.Cells.Item(2, 12).Value = 1105.2
$ws.Cells.Item(2, 14).Value = -1331.2
$ws.Cells.Item(3, 8).Value = 7662.4443
$ws.Cells.Item(3, 9).Value = 7423.143
$ws.Cells.Item(3, 10).Value = 8500
$ws.Cells.Item(3, 11).Value = 22269.429
$ws.Cells.Item(3, 12).Value = 25500
$ws.Cells.Item(3, 13).Value = -22157.429
$ws.Cells.Item(3, 14).Value = -25724
$ws.Cells.Item(9, 8).Value = 183367120
$ws.Cells.Item(9, 10).Value = 250050600
$ws.Cells.Item(9, 12).Value = 750151800
$ws.Cells.Item(9, 14).Value = -750152248
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).ClearContents()
$ws.Cells.Item(46, 14).ClearContents()
$ws.Cells.Item(58, 8).Value = 46757.855
$ws.Cells.Item(58, 9).Value = 301.66666
$ws.Cells.Item(58, 10).Value = 81600
$ws.Cells.Item(58, 11).Value = 904.9999799999999
$ws.Cells.Item(58, 12).Value = 244800
$ws.Cells.Item(58, 13).Value = -776.9999799999999
$ws.Cells.Item(58, 14).Value = -245056
$ws.Cells.Item(64, 8).Value = 2000
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 2000
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 6000
$ws.Cells.Item(64, 13).ClearContents()
$ws.Cells.Item(64, 14).Value = -6540
$ws.Cells.Item(67, 8).Value = 2000
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 2000
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 12).Value = 6000
$ws.Cells.Item(67, 13).ClearContents()
$ws.Cells.Item(67, 14).Value = -7872
$ws.Cells.Item(81, 8).Value = 1245.6
$ws.Cells.Item(81, 9).Value = 407.66666
$ws.Cells.Item(81, 10).Value = 2502.5
$ws.Cells.Item(81, 11).Value = 1222.99998
$ws.Cells.Item(81, 12).Value = 7507.5
$ws.Cells.Item(81, 13).Value = -99.99998000000005
$ws.Cells.Item(81, 14).Value = -9753.5
$ws.Cells.Item(84, 8).Value = 1245.6
$ws.Cells.Item(84, 9).Value = 407.66666
$ws.Cells.Item(84, 10).Value = 2502.5
$ws.Cells.Item(84, 11).Value = 3668.99994
$ws.Cells.Item(84, 12).Value = 22522.5
$ws.Cells.Item(84, 13).Value = 1947.00006
$ws.Cells.Item(84, 14).Value = -33754.5
$ws.Cells.Item(106, 8).Value = 2009.6666
$ws.Cells.Item(106, 10).Value = 2009.6666
$ws.Cells.Item(106, 12).Value = 6028.9998
$ws.Cells.Item(106, 14).Value = -7920.9998
$ws.Cells.Item(109, 8).Value = 2613.9666
$ws.Cells.Item(109, 9).Value = 617.38464
$ws.Cells.Item(109, 10).Value = 4140.7646
$ws.Cells.Item(109, 11).Value = 1852.15392
$ws.Cells.Item(109, 12).Value = 12422.2938
$ws.Cells.Item(109, 13).Value = -812.15392
$ws.Cells.Item(109, 14).Value = -14502.2938
$ws.Cells.Item(120, 8).Value = 14275
$ws.Cells.Item(120, 9).Value = 6400
$ws.Cells.Item(120, 11).Value = 19200
$ws.Cells.Item(120, 13).Value = -14362
$ws.Cells.Item(121, 8).Value = 18519706
$ws.Cells.Item(121, 9).Value = 700
$ws.Cells.Item(121, 10).Value = 22223508
$ws.Cells.Item(121, 11).Value = 2100
$ws.Cells.Item(121, 12).Value = 66670524
$ws.Cells.Item(121, 13).Value = -790
$ws.Cells.Item(121, 14).Value = -66673144
$ws.Cells.Item(122, 8).Value = 35721804
$ws.Cells.Item(122, 10).Value = 50999.5
$ws.Cells.Item(122, 12).Value = 458995.5
$ws.Cells.Item(122, 14).Value = -463895.5
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(124, 13).ClearContents()
$ws.Cells.Item(124, 14).ClearContents()
$ws.Cells.Item(125, 8).Value = 3900
$ws.Cells.Item(125, 9).Value = 3900
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 11).Value = 11700
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 13).Value = -6780
$ws.Cells.Item(125, 14).ClearContents()
$ws.Cells.Item(134, 8).Value = 22728868
$ws.Cells.Item(134, 9).Value = 26316864
$ws.Cells.Item(134, 10).Value = 4899.6665
$ws.Cells.Item(134, 11).Value = 78950592
$ws.Cells.Item(134, 12).Value = 14698.9995
$ws.Cells.Item(134, 13).Value = -78945522
$ws.Cells.Item(134, 14).Value = -24838.9995
$ws.Cells.Item(139, 8).Value = 7737984.5
$ws.Cells.Item(139, 9).Value = 11364196
$ws.Cells.Item(139, 10).Value = 485562.62
$ws.Cells.Item(139, 11).Value = 34092588
$ws.Cells.Item(139, 12).Value = 1456687.86
$ws.Cells.Item(139, 13).Value = -34087448
$ws.Cells.Item(139, 14).Value = -1466967.86
$ws.Cells.Item(140, 8).Value = 18751852
$ws.Cells.Item(140, 9).Value = 45000856
$ws.Cells.Item(140, 10).Value = 2564.2144
$ws.Cells.Item(140, 11).Value = 135002568
$ws.Cells.Item(140, 12).Value = 7692.6432
$ws.Cells.Item(140, 13).Value = -134997388
$ws.Cells.Item(140, 14).Value = -18052.6432
$ws.Cells.Item(141, 8).Value = 50002560
$ws.Cells.Item(141, 9).Value = 58825624
$ws.Cells.Item(141, 10).Value = 5200
$ws.Cells.Item(141, 11).Value = 176476872
$ws.Cells.Item(141, 12).Value = 15600
$ws.Cells.Item(141, 13).Value = -176471692
$ws.Cells.Item(141, 14).Value = -25960
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 41.42857
$ws.Cells.Item(2, 9).Value = 15
$ws.Cells.Item(2, 10).Value = 76.666664
$ws.Cells.Item(2, 11).Value = 15
$ws.Cells.Item(2, 12).Value = 76.666664
$ws.Cells.Item(2, 13).Value = 98
$ws.Cells.Item(2, 14).Value = -302.666664
$ws.Cells.Item(70, 8).Value = 3995.0908
$ws.Cells.Item(70, 10).Value = 4062.8572
$ws.Cells.Item(70, 12).Value = 4062.8572
$ws.Cells.Item(70, 14).Value = -4602.8572
$ws.Cells.Item(73, 8).Value = 3995.0908
$ws.Cells.Item(73, 10).Value = 4062.8572
$ws.Cells.Item(73, 12).Value = 4062.8572
$ws.Cells.Item(73, 14).Value = -5934.8572
$ws.Cells.Item(102, 8).Value = 1406.52
$ws.Cells.Item(102, 9).Value = 1182.7222
$ws.Cells.Item(102, 11).Value = 1182.7222
$ws.Cells.Item(102, 13).Value = 439.2778000000001
$ws.Cells.Item(113, 8).Value = 1382.9166
$ws.Cells.Item(113, 10).Value = 1649
$ws.Cells.Item(113, 12).Value = 1649
$ws.Cells.Item(113, 14).Value = -5989
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 63026360
$ws.Cells.Item(16, 9).Value = 5495130
$ws.Cells.Item(16, 10).Value = 250002850
$ws.Cells.Item(16, 11).Value = 5495130
$ws.Cells.Item(16, 12).Value = 250002850
$ws.Cells.Item(16, 13).Value = -5494960
$ws.Cells.Item(16, 14).Value = -250003190
